$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 379-392 ---
$ws = $wb.Worksheets.Item("PIR")
$pirRows = @(
  @('2026-02-04','14:32:39','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:32:42','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:32:44','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:32:50','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:32:55','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:33:00','14:00','Bathroom','Motion Detected','Active'),
  @('2026-02-04','14:33:06','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:33:10','14:00','Bathroom','Motion Detected','Active'),
  @('2026-02-04','14:33:16','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:33:17','14:00','Bathroom','Motion Detected','Active'),
  @('2026-02-04','14:33:25','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:33:27','14:00','Bathroom','Motion Detected','Active'),
  @('2026-02-04','14:33:35','14:00','Bathroom','No Motion','Inactive'),
  @('2026-02-04','14:33:36','14:00','Bathroom','Motion Detected','Active')
)
$startRow = 379
$r = $startRow
foreach ($row in $pirRows) {
  $ws.Cells.Item($r, 1).Value = "'" + $row[0]
  $ws.Cells.Item($r, 2).Value = "'" + $row[1]
  $ws.Cells.Item($r, 3).Value = "'" + $row[2]
  $ws.Cells.Item($r, 4).Value = "'" + $row[3]
  $ws.Cells.Item($r, 5).Value = "'" + $row[4]
  $ws.Cells.Item($r, 6).Value = "'" + $row[5]
  $r = $r + 1
}

# --- Humidity sheet: append rows 308-319 ---
$ws = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
  @('2026-02-04','14:32:37','14:00','Bathroom','78.8%','Active'),
  @('2026-02-04','14:32:40','14:00','Bathroom','79.9%','Active'),
  @('2026-02-04','14:32:48','14:00','Bathroom','79.8%','Active'),
  @('2026-02-04','14:32:53','14:00','Bathroom','78.7%','Active'),
  @('2026-02-04','14:32:58','14:00','Bathroom','79.6%','Active'),
  @('2026-02-04','14:33:03','14:00','Bathroom','78.7%','Active'),
  @('2026-02-04','14:33:08','14:00','Bathroom','79.6%','Active'),
  @('2026-02-04','14:33:13','14:00','Bathroom','78.6%','Active'),
  @('2026-02-04','14:33:18','14:00','Bathroom','79.6%','Active'),
  @('2026-02-04','14:33:23','14:00','Bathroom','78.7%','Active'),
  @('2026-02-04','14:33:28','14:00','Bathroom','79.6%','Active'),
  @('2026-02-04','14:33:33','14:00','Bathroom','78.8%','Active')
)
$startRow = 308
$r = $startRow
foreach ($row in $humidityRows) {
  $ws.Cells.Item($r, 1).Value = "'" + $row[0]
  $ws.Cells.Item($r, 2).Value = "'" + $row[1]
  $ws.Cells.Item($r, 3).Value = "'" + $row[2]
  $ws.Cells.Item($r, 4).Value = "'" + $row[3]
  $ws.Cells.Item($r, 5).Value = "'" + $row[4]
  $ws.Cells.Item($r, 6).Value = "'" + $row[5]
  $r = $r + 1
}

# --- Temperature sheet: append rows 308-319 ---
$ws = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
  @('2026-02-04','14:32:38','14:00','Bathroom','24.2C','Active'),
  @('2026-02-04','14:32:41','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:32:49','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:32:54','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:32:59','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:33:04','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:33:09','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:33:14','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:33:19','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:33:24','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:33:29','14:00','Bathroom','24.3C','Active'),
  @('2026-02-04','14:33:34','14:00','Bathroom','24.3C','Active')
)
$startRow = 308
$r = $startRow
foreach ($row in $temperatureRows) {
  $ws.Cells.Item($r, 1).Value = "'" + $row[0]
  $ws.Cells.Item($r, 2).Value = "'" + $row[1]
  $ws.Cells.Item($r, 3).Value = "'" + $row[2]
  $ws.Cells.Item($r, 4).Value = "'" + $row[3]
  $ws.Cells.Item($r, 5).Value = "'" + $row[4]
  $ws.Cells.Item($r, 6).Value = "'" + $row[5]
  $r = $r + 1
}

